# Add the new "2022-Q1" sheet with fund-holding detail data, positioned
# immediately before the "总计" (total) summary sheet. The sheet reference
# used to anchor the insertion position is not reused afterwards -- this
# COM shim invalidates previously-held Worksheet object handles once
# Worksheets.Add() runs, so every later lookup re-resolves by name.
$wb = $excel.ActiveWorkbook
$q1 = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$q1.Name = "2022-Q1"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($c = 2; $c -le 8; $c++) {
    $cell = $q1.Cells.Item(1, $c)
    $cell.Value = $headers[$c - 2]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$dataRows = @(
    @('001445', '华安国企改革主题灵活配置混合', '47.50', '87.35', '3.91', '1.8572', 6),
    @('002803', '东方红沪港深灵活配置混合', '41.94', '83.46', '3.59', '1.5056', 7),
    @('010341', '招商产业精选股票A', '30.60', '81.87', '4.90', '1.4994', 6),
    @('000746', '招商行业精选股票', '26.52', '83.71', '4.88', '1.2942', 7),
    @('161706', '招商优质成长混合 (LOF)', '16.57', '80.28', '4.84', '0.8020', 7),
    @('012835', '招商景气精选股票型证券投资基金A', '11.18', '83.90', '4.95', '0.5534', 6),
    @('011568', '鹏华产业升级混合A', '21.41', '68.03', '2.12', '0.4539', 7),
    @('169102', '东方红睿阳三年定期开放灵活配置混合', '8.68', '81.40', '3.59', '0.3116', 9),
    @('009601', '招商科技动力3个月滚动持有股票A', '5.07', '82.64', '4.91', '0.2489', 6),
    @('012093', '鹏华创新升级混合型证券投资基金A', '6.58', '64.31', '2.65', '0.1744', 6),
    @('012836', '招商景气精选股票型证券投资基金C', '3.36', '83.90', '4.95', '0.1663', 6),
    @('008261', '招商研究优选股票A', '2.86', '83.26', '4.92', '0.1407', 7),
    @('004206', '华商元亨灵活配置混合', '5.64', '29.94', '1.88', '0.1060', 3),
    @('519625', '银河君盛灵活配置混合A', '4.91', '20.05', '2.00', '0.0982', 2),
    @('770001', '德邦优化灵活配置混合', '2.49', '86.80', '3.24', '0.0807', 5),
    @('700001', '平安行业先锋混合', '2.32', '91.67', '3.05', '0.0708', 6),
    @('011071', '鹏华安悦一年持有期混合A', '9.16', '21.81', '0.69', '0.0632', 4),
    @('003598', '华商润丰灵活配置混合A', '3.17', '37.86', '1.94', '0.0615', 4),
    @('007509', '华商润丰灵活配置混合C', '3.09', '37.86', '1.94', '0.0599', 4),
    @('006167', '德邦乐享生活混合A', '1.34', '88.13', '3.56', '0.0477', 10),
    @('519626', '银河君盛灵活配置混合C', '2.33', '20.05', '2.00', '0.0466', 2),
    @('001067', '鹏华弘盛灵活配置混合A', '6.22', '20.46', '0.69', '0.0429', 4),
    @('009169', '湘财长兴灵活配置混合A', '1.16', '85.40', '3.18', '0.0369', 6),
    @('008840', '德邦大消费混合A', '0.95', '89.43', '3.44', '0.0327', 9),
    @('009602', '招商科技动力3个月滚动持有股票C', '0.54', '82.64', '4.91', '0.0265', 6),
    @('009232', '鹏华安惠混合A', '3.74', '21.93', '0.64', '0.0239', 6),
    @('010257', '天弘多利一年定期开放混合', '2.64', '26.41', '0.87', '0.0230', 5),
    @('006168', '德邦乐享生活混合C', '0.42', '88.13', '3.56', '0.0150', 10),
    @('009170', '湘财长兴灵活配置混合C', '0.46', '85.40', '3.18', '0.0146', 6),
    @('008841', '德邦大消费混合C', '0.34', '89.43', '3.44', '0.0117', 9),
    @('001448', '华商双翼平衡混合', '0.38', '39.74', '2.85', '0.0108', 3),
    @('001664', '平安鑫安混合A', '0.86', '29.46', '1.11', '0.0095', 3),
    @('006225', '人保量化基本面混合A', '0.63', '88.00', '1.50', '0.0094', 10),
    @('008135', '华宸未来价值先锋混合', '0.20', '86.99', '3.66', '0.0073', 10),
    @('011569', '鹏华产业升级混合C', '0.34', '68.03', '2.12', '0.0072', 7),
    @('570006', '诺德中小盘混合', '0.20', '90.18', '2.81', '0.0056', 4),
    @('007049', '平安鑫安混合E', '0.50', '29.46', '1.11', '0.0056', 3),
    @('012094', '鹏华创新升级混合型证券投资基金C', '0.19', '64.31', '2.65', '0.0050', 6),
    @('003416', '招商财经大数据策略股票A', '0.09', '82.37', '4.94', '0.0044', 6),
    @('003626', '平安鑫利灵活配置混合A', '0.41', '27.68', '0.98', '0.0040', 5),
    @('001380', '鹏华弘盛灵活配置混合C', '0.53', '20.46', '0.69', '0.0037', 4),
    @('006433', '平安鑫利灵活配置混合C', '0.25', '27.68', '0.98', '0.0024', 5),
    @('009233', '鹏华安惠混合C', '0.32', '21.93', '0.64', '0.0020', 6),
    @('006226', '人保量化基本面混合C', '0.04', '88.00', '1.50', '0.0006', 10),
    @('005493', '鑫元价值精选灵活配置混合A', '0.07', '21.81', '0.69', '0.0005', 4),
    @('007952', '招商财经大数据策略股票C', '0.01', '82.37', '4.94', '0.0005', 6),
    @('001665', '平安鑫安混合C', '0.02', '29.46', '1.11', '0.0002', 3)
)

$r = 2
foreach ($row in $dataRows) {
    $a = $q1.Cells.Item($r, 1)
    $a.Value = ($r - 2)
    $a.Font.Bold = $true
    $a.HorizontalAlignment = -4108
    $a.VerticalAlignment = -4160
    $a.Borders.LineStyle = 1

    for ($c = 2; $c -le 7; $c++) {
        $cell = $q1.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$c - 2]
        $cell.ClearFormats()
    }

    $q1.Cells.Item($r, 8).Value = $row[6]

    $r++
}

# Insert the new "2022-Q1" summary row at the top of the data in "总计",
# pushing the existing quarters down by one row. Re-fetch the sheet by
# name (see note above) rather than reusing the handle from the Add() call.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$a2 = $totalSheet.Cells.Item(2, 1)
$a2.Value = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 47
$totalSheet.Cells.Item(2, 4).Value = 9.949999999999999

Write-Output "edit complete"
